$d = $word.ActiveDocument

# The "Medicine for Claxter" scene has two almost-identical opening lines:
#   "Ornev: [random] Claxter! How are you doing today?"
#   "Ornev: [random] Hello, how are you?"
# The second one is an old leftover/duplicate line that should be removed
# entirely, and the now-redundant "[random]" annotation should be dropped
# from the line that remains.

# 1. Delete the duplicate paragraph "Ornev: [random] Hello, how are you?"
#    in its entirety, including its paragraph mark.
$dup = $d.Content
$foundDup = $dup.Find.Execute("Ornev: [random] Hello, how are you?", $true,
                               $false, $false, $false, $false, $true, 1,
                               $false, "", 0)
if ($foundDup) {
    $dupPara = $d.Range($dup.Start, $dup.End + 1)
    $dupPara.Delete()
}

# 2. Remove the now-redundant "[random] " annotation (and its trailing
#    space run) from the remaining paragraph, which currently reads:
#    "Ornev: [random] Claxter! How are you doing today?"
#    This leaves the original bold space run that sat between "Ornev:"
#    and "[random]" in place, so the result reads
#    "Ornev: Claxter! How are you doing today?" with a single space.
$annot = $d.Content
$foundAnnot = $annot.Find.Execute("[random] ", $true, $false, $false,
                                   $false, $false, $true, 1, $false, "", 0)
if ($foundAnnot) {
    $annot.Delete()
}
